# Updating burndown sprint1 chart
#
# Sprint1 backlog: fill in the Week 3 / Week 4 "amount remaining" columns
# (F and G) for the tasks that had just been completed/updated, which
# ripples through the Estimate Totals row (34) and, via Sheet1!$C$34:$G$34,
# into the burndown chart's data series.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9  - "Testing model/ viewmodel": work wrapped up, 0 remaining
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0

# Row 15 - "Testing viewmodel of subpages": wrapped up in week 4
$ws.Range("G15").Value = 0

# Row 17 - "UI for ability to select new game user owns": progress this week
$ws.Range("F17").Value = 3
$ws.Range("G17").Value = 2

# Row 18 - "UI to display details from game, use image here": no progress yet
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 4

# Move the active selection to reflect where work left off
$ws.Range("G17").Select()

$wb.Save()
